$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutomationModules")

# TM003 row: rename module and flip Execute flag to Y
$ws.Range("B4").Value = "Login_into_EMR"
$ws.Range("D4").Value = "Y"

# TM002 row: fill in Module Name and EXCEL_KEY with "Registration"
$ws.Range("B3").Value = "Registration"
$ws.Range("C3").Value = "Registration"

# Update the selected cell on the active sheet
$ws.Range("C3").Select()
